$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A163").Value = "B018"
$ws.Range("B163").Value = "Macaíba"
$ws.Range("C163").Value = "NE"

$ws.Range("A164").Value = "B021"
$ws.Range("B164").Value = "Mossóro"
$ws.Range("C164").Value = "NO"

$ws.Range("A165").Value = "B001"
$ws.Range("B165").Value = "São Gonçalo Amarante"
$ws.Range("C165").Value = "NE"

$ws.Range("H163").Select()
$excel.ActiveWindow.ScrollRow = 151
